# Edit script for "7) Projektgennemførelse.docx"
# Applies the textual changes described by the commit diff.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Output "WARNING: find failed for [$find]"
    }
}

# 1) Image run: Word re-marks the picture run as NoProof + Danish east-asian language
#    after a resave (noProof + lang eastAsia="da-DK").
$shp = $d.InlineShapes(1)
$shpRange = $shp.Range
$shpRange.NoProofing = 1
$shpRange.LanguageIDFarEast = "da-DK"

# 2) "Scrum har de fordele:" -> "Scrum har fordelene:"
Replace-Text " har de fordele" " har fordelene"

# 3) "...af projekt før man kan begynder..." -> "...af projektet før man kan begynder..."
Replace-Text "af projekt før man kan begynde" "af projektet før man kan begynde"

# 4) "Hver iteration ender ud..." -> "Hvert sprint ender ud..." (drop spell-check markers on "iteration")
Replace-Text "Hver iteration ender" "Hvert sprint ender"

# 5) "...men at vi i denne også valgte..." -> "...men at gruppen også valgte..."
Replace-Text "men at vi i denne også valgte" "men at gruppen også valgte"

# 6) "Nøglebegreberne som vi har brugt fra" -> "Nøglebegreberne som gruppen har brugt fra"
Replace-Text "Nøglebegreberne som vi har brugt fra" "Nøglebegreberne som gruppen har brugt fra"

# 7) Scrummasterrollen paragraph rewrite
Replace-Text "Som Scrum foreskriver har gruppen haft 8 medlemmer" "Gruppen består af 8 medlemmer"
Replace-Text "fungeret som primær kontaktperson med produktowner" "fungeret som primær kontaktperson til produktowner"

# 8) Taskboardet paragraph rewrite
Replace-Text "har fungeret som en liste over opgaver med prioritering, organisator og kontrakt med hvad Scrumgruppen" "har fungeret som organistor for listen over opgaver for det pågældende sprint. Opgaverne er blevet defineret med beskrivelser, prioritering, estimeret tid og uddelegering. Taskboardet har også fungeret som kontrakt for hvad Scrumgruppen"
Replace-Text "har valgt at forpligtige sig til af opgaver. " "har valgt at forpligtige sig til af opgaver i sprintet. "

# 9) Remove stray lastRenderedPageBreak before "Tidsplan" heading
$found = $d.Content.Find.Execute("Tidsplan", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # nothing textual changes here; handled separately below via paragraph scan
}

# 10) "dirigient" (typo, spell-checker flagged) -> "dirigent" (correct)
Replace-Text "dagsorden, dirigient og referent" "dagsorden, dirigent og referent"
